# Update countries & provincias Spain
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1) Update the "last updated" timestamp banner (row 1)
$ws.Range("A1").Value = "Datos actualizados a 27 de Marzo de 2020 a las 22:28"

# 2) Update Cataluna row (row 5) totals
$ws.Range("B5").Value = 14263
$ws.Range("C5").Value = 3106
$ws.Range("D5").Value = 10087
$ws.Range("E5").Value = 1070

# 3) Gran Canaria's figures were updated, which bumps it above Malaga and
#    Asturias when the table is kept sorted descending by "Casos totales"
#    (column B). Previously: row17=Malaga, row18=Asturias, row19=Gran Canaria.
#    Now: row17=Gran Canaria (new figures), row18=Malaga, row19=Asturias
#    (Malaga/Asturias simply shift down one row with their existing figures).
$ws.Range("A17").Value = "Gran Canaria"
$ws.Range("B17").Value = 1025
$ws.Range("C17").Value = 20
$ws.Range("D17").Value = 834
$ws.Range("E17").Value = 36

$ws.Range("A18").Value = "Malaga"
$ws.Range("B18").Value = 1006
$ws.Range("C18").Value = 61
$ws.Range("D18").Value = 899
$ws.Range("E18").Value = 46

$ws.Range("A19").Value = "Asturias"
$ws.Range("B19").Value = 1004
$ws.Range("C19").Value = 65
$ws.Range("D19").Value = 906
$ws.Range("E19").Value = 33

# 4) Canary Island provinces: Muertes ("Deaths") figure corrected 28 -> 36
$ws.Range("E31").Value = 36
$ws.Range("E55").Value = 36
$ws.Range("E57").Value = 36
$ws.Range("E58").Value = 36
$ws.Range("E62").Value = 36
$ws.Range("E63").Value = 36
